# clarke_exploring_archives.docx -- "change pdfs and docs"
#
# 1. Heading "Peer-Reviewed Assessment" -> "Peer-Reviewed Assignment"
# 2. Tidy "..., Saidiya Hartman warns that" so it is one clean run of text
#    (drops the spell-check proofErr markers that bracketed "Saidiya")

$d = $word.ActiveDocument

# 1. Assessment -> Assignment
$d.Content.Find.Execute("Assessment", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Assignment", 2) | Out-Null

# 2. Re-write the "Saidiya Hartman" phrase as a single run of plain text
$d.Content.Find.Execute(", Saidiya Hartman warns that", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", Saidiya Hartman warns that", 2) | Out-Null
